# Applies the "Add files via upload" trading-journal update:
# - extends the existing note in intraday!E81
# - appends new dated rows to intraday (sheet1), learnings (sheet3) and market_sR (sheet7)
# - reuses the existing date style (m/d/yyyy) via copy/paste-format so no duplicate number formats are created
# - applies three new solid fills (red / light blue / green) introduced by the edit

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("intraday")
$sheet3 = $wb.Worksheets.Item("learnings")
$sheet7 = $wb.Worksheets.Item("market_sR")

# --- helper: stamp a date value onto a cell reusing the workbooks existing
#     "m/d/yyyy" date style (sheet1!A81) instead of letting Excel synthesize a new one
function Set-DateCell {
    param($range, $serial)
    $sheet1.Range("A81").Copy() | Out-Null
    $range.PasteSpecial(-4122) | Out-Null
    $range.Value2 = $serial
    $excel.CutCopyMode = 0
}

# --- intraday!E81: append the missed-profit postscript to the existing note
$sheet1.Range("E81").Value2 = $sheet1.Range("E81").Value2 + "/if I would have waited today I would have made 18k plus profit in single lot"

# --- intraday row 82
Set-DateCell $sheet1.Range("A82") 45664
$sheet1.Range("B82").Value2 = 1317
$sheet1.Range("C82").Value2 = " "
$sheet1.Range("D82").Value2 = 3
$sheet1.Range("E82").Value2 = "don't trade emotionally Man always plot some resistance and support/// always wait for trade  I knew there was down trade but my entry was abit fast"

# --- intraday row 83
Set-DateCell $sheet1.Range("A83") 45665
$sheet1.Range("C83").Value2 = 1365
$sheet1.Range("D83").Value2 = 2
$sheet1.Range("E83").Value2 = "Because of groww SL mistake I booked abit big loss today and market seems not in my favour/// Market tried to move up but there were three candles in 5 min frame where they didn't break that so my entry should be threre"

# --- intraday row 84
Set-DateCell $sheet1.Range("A84") 45666
$sheet1.Range("C84").Value2 = 1438
$sheet1.Range("D84").Value2 = 3
$sheet1.Range("E84").Value2 = "I broke my rule that is if I would have been patinece then it would be good/ learn to be patience always in trading and wait for the setup to break"

# --- intraday row 85: blank separator cell with a red fill
$sheet1.Range("B85").Interior.Color = 255

# --- intraday row 86: section note
$sheet1.Range("B86").Value2 = " Two Trade from 01/09/2025 for at least 6months and see the performance// and I will follow my rule and I will earn as well"

# --- intraday row 87
Set-DateCell $sheet1.Range("A87") 45667
$sheet1.Range("B87").Value2 = 556
$sheet1.Range("D87").Value2 = 2
$sheet1.Range("E87").Value2 = "first trade was an emotional trade and second trade was perfect and I booked 47 poiunt which was good  but always track your trade don’t panic"

# --- learnings row 9 (new unrelated entry jotted down amid the intraday rows)
Set-DateCell $sheet3.Range("A9") 45931
$sheet3.Range("B9").Value2 = "after taking a trade always try to engage with trade not with fear "

# --- intraday row 88
Set-DateCell $sheet1.Range("A88") 45670
$sheet1.Range("B88").Value2 = 200
$sheet1.Range("C88").Value2 = " "
$sheet1.Range("D88").Value2 = 3
$sheet1.Range("E88").Value2 = "first rade was ok if that was good I would have booked good profit but second trade was based on other opinion try to lean and earn in market"

# --- intraday row 89: B carries a note with a light-blue fill, E a separate reason
$sheet1.Range("B89").Value2 = "I prmosed that I will trade only 1 and 1 but one in morning and one in eveing after knowing market trend"
$sheet1.Range("B89").Interior.Color = 15773696
$sheet1.Range("E89").Value2 = "Never ever be in hurry to take trade as per your thinking alwasys wait for setup and support and resistence to break "

# --- market_sR row 3
Set-DateCell $sheet7.Range("A3") 45670
$sheet7.Range("B3").Value2 = "market opened gap down and it moved abit up and again fell down later "

# --- market_sR row 4
Set-DateCell $sheet7.Range("A4") 45671
$sheet7.Range("B4").Value2 = "market opened gap up and moved abit up and down again "

# --- intraday row 90
Set-DateCell $sheet1.Range("A90") 45671
$sheet1.Range("B90").Value2 = 814
$sheet1.Range("D90").Value2 = 1
$sheet1.Range("E90").Value2 = "it was perfect setup as per my thinking"

# --- intraday row 91: E (reason) is filled before B (note, green fill)
Set-DateCell $sheet1.Range("A91") 45672
$sheet1.Range("C91").Value2 = " "
$sheet1.Range("D91").Value2 = 3
$sheet1.Range("E91").Value2 = "both trade were against my trade so entry was ok but market didn't move as per my logic///but again the previous logic which I marked worked perfectly but I didn’t trade // always my last trade goes well so I need to wait again  but learn to wait abit more seeing the price action"
$sheet1.Range("B91").Value2 = "I placed SL perfectly but I always missed to book big profit so try to hold trade abit longer                                                    165"
$sheet1.Range("B91").Interior.Color = 5287936

# --- market_sR row 5
Set-DateCell $sheet7.Range("A5") 45672
$sheet7.Range("B5").Value2 = "maket opengap up and went abit up and moved down again "

# --- market_sR row 6
Set-DateCell $sheet7.Range("A6") 45673
$sheet7.Range("B6").Value2 = "market open gap up and went abit up and moved down again "

# --- intraday row 92 (no date, closes the block)
$sheet1.Range("B92").Value2 = 1262
$sheet1.Range("D92").Value2 = 2
$sheet1.Range("E92").Value2 = "patience is the key and level of support and resistence is also important"

